$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$unsafeCells = @("D5","D6","D7","D11","D13","D15","D16","D19","D22","D23","D24","D26","D27","D28","D29","D30","D35","D38","D39","D40","D41","D45")
foreach ($ref in $unsafeCells) { $ws.Range($ref).NumberFormat = "@" }

$ws.Range("D2").Value = "70.224.61"
$ws.Range("E2").Value = "  +0.33%  "

$ws.Range("D3").Value = "3.599.95"
$ws.Range("E3").Value = "  +1.50%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").Value = "603.89"
$ws.Range("E5").Value = "  +0.10%  "

$ws.Range("D6").Value = "195.50"
$ws.Range("E6").Value = "  -0.73%  "

$ws.Range("D7").Value = "0.626"
$ws.Range("E7").Value = "  -0.14%  "

$ws.Range("E8").Value = "  +0.05%  "

$ws.Range("E9").Value = "  -1.97%  "

$ws.Range("E10").Value = "  -1.64%  "

$ws.Range("D11").Value = "53.72"
$ws.Range("E11").Value = "  -0.75%  "

$ws.Range("E12").Value = "  +0.03%  "

$ws.Range("D13").Value = "9.54"
$ws.Range("E13").Value = "  -0.38%  "

$ws.Range("D14").Value = "4.174.55"
$ws.Range("E14").Value = "  +1.84%  "

$ws.Range("D15").Value = "13.05"
$ws.Range("E15").Value = "  +2.70%  "

$ws.Range("D16").Value = "595.74"
$ws.Range("E16").Value = "  -0.92%  "

$ws.Range("D17").Value = "70.356.22"
$ws.Range("E17").Value = "  +0.31%  "

$ws.Range("D18").Value = "3.602.65"
$ws.Range("E18").Value = "  +1.73%  "

$ws.Range("D19").Value = "19.03"
$ws.Range("E19").Value = "  -0.81%  "

$ws.Range("E20").Value = "  +1.38%  "

$ws.Range("E21").Value = "  -0.31%  "

$ws.Range("D22").Value = "17.81"
$ws.Range("E22").Value = "  -1.84%  "

$ws.Range("D23").Value = "5.19"
$ws.Range("E23").Value = "  -2.83%  "

$ws.Range("D24").Value = "102.06"
$ws.Range("E24").Value = "  -1.52%  "

$ws.Range("E25").Value = "  -0.34%  "

$ws.Range("D26").Value = "3.02"
$ws.Range("E26").Value = "  -3.24%  "

$ws.Range("D27").Value = "10.75"
$ws.Range("E27").Value = "  -2.00%  "

$ws.Range("D28").Value = "9.63"
$ws.Range("E28").Value = "  -0.73%  "

$ws.Range("D29").Value = "33.79"
$ws.Range("E29").Value = "  +0.43%  "

$ws.Range("D30").Value = "4.77"
$ws.Range("E30").Value = "  +6.12%  "

$ws.Range("E31").Value = "  +0.27%  "

$ws.Range("E32").Value = "  -3.85%  "

$ws.Range("E33").Value = "  +0.82%  "

$ws.Range("D34").Value = "0.0₃0902"
$ws.Range("E34").Value = "  +8.44%  "

$ws.Range("D35").Value = "63.18"
$ws.Range("E35").Value = "  -0.45%  "

$ws.Range("D36").Value = "3.897.51"
$ws.Range("E36").Value = "  +4.23%  "

$ws.Range("E37").Value = "  -0.36%  "

$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "523.80"
$ws.Range("E38").Value = "  +4.87%  "

$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").Value = "1.00"
$ws.Range("E39").Value = "  -0.07%  "

$ws.Range("D40").Value = "36.86"
$ws.Range("E40").Value = "  -0.16%  "

$ws.Range("D41").Value = "0.390"
$ws.Range("E41").Value = "  -1.39%  "

$ws.Range("E42").Value = "  -2.83%  "

$ws.Range("E43").Value = "  -2.06%  "

$ws.Range("E44").Value = "  -0.80%  "

$ws.Range("D45").Value = "3.39"
$ws.Range("E45").Value = "  +2.10%  "

$ws.Range("E46").Value = "  +0.46%  "

$ws.Range("E47").Value = "  +0.01%  "

$ws.Range("E48").Value = "  -0.85%  "

$ws.Range("E49").Value = "  -0.18%  "

$ws.Range("E50").Value = "  +2.09%  "

$ws.Range("E51").Value = "  +0.50%  "

foreach ($ref in $unsafeCells) { $ws.Range($ref).ClearFormats() }
